$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 value from 6 to 8
$ws.Range("A2").Value = 8

# Add new row 8 data
$ws.Range("B8").Value = 2500
$ws.Range("C8").Value = 7
$ws.Range("D8").Value = 8

# Add new row 9 data
$ws.Range("B9").Value = 5000
$ws.Range("C9").Value = 8
$ws.Range("D9").Value = 7

# Update the active cell selection to F8
$ws.Range("F8").Select()
